# Auto-generated: applies cryptos.xlsx cell-value updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.681.86"
$ws.Range("E2").Value = "  -4.06%  "
$ws.Range("D3").Value = "2.370.52"
$ws.Range("E3").Value = "  +4.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "292.77"
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.96"
$ws.Range("E6").Value = "  -5.62%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  -1.50%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.496"
$ws.Range("E9").Value = "  -3.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.83"
$ws.Range("E10").Value = "  -5.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0773"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.89"
$ws.Range("E13").Value = "  -4.23%  "
$ws.Range("D14").Value = "2.719.92"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("D15").Value = "2.361.69"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.85"
$ws.Range("E16").Value = "  +1.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.817"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "44.653.11"
$ws.Range("E18").Value = "  -4.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.29"
$ws.Range("E19").Value = "  -4.57%  "
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.76"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.12"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.74"
$ws.Range("E24").Value = "  -3.88%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.20"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.83"
$ws.Range("E28").Value = "  -9.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.48"
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.81"
$ws.Range("E30").Value = "  +16.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.81"
$ws.Range("E31").Value = "  +5.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "147.07"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.68"
$ws.Range("E33").Value = "  -4.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0753"
$ws.Range("E35").Value = "  -3.09%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.94"
$ws.Range("E36").Value = "  +11.78%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.110"
$ws.Range("E37").Value = "  -4.16%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.70"
$ws.Range("E39").Value = "  -7.49%  "
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.14"
$ws.Range("E42").Value = "  -4.09%  "
$ws.Range("D43").Value = "1.927.43"
$ws.Range("E43").Value = "  +6.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "89.46"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.70"
$ws.Range("E46").Value = "  -12.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.46"
$ws.Range("E47").Value = "  +8.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.02"
$ws.Range("E48").Value = "  +15.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "99.44"
$ws.Range("E49").Value = "  +5.29%  "
$ws.Range("D50").Value = "2.591.16"
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.181"
$ws.Range("E51").Value = "  -4.96%  "
